$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the old "_GoBack" bookmark (currently sitting at the end
#    of the ". I need to discover what finger she lands on for:"
#    paragraph). It will be re-created further down, in the new
#    closing paragraph of the answer to question 3.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $oldBm = $d.Bookmarks.Item("_GoBack")
    $oldBm.Delete()
}

# ------------------------------------------------------------------
# 2. Locate the "Count of 1000" paragraph (the last line of the
#    "what finger does she land on" list) so we can insert the new
#    material for question 3 right after it.
# ------------------------------------------------------------------
$findRange = $d.Content
$findRange.Find.Execute("Count of 1000", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$countPara = $findRange.Paragraphs.First
$insertAt = $d.Range($countPara.Range.Start, $countPara.Range.End - 1)
$insertAt.Collapse(0)

# ------------------------------------------------------------------
# 3. Insert the new content: a blank line, a blank (but formatted)
#    paragraph, the "3. Identify potential solutions" heading and
#    the answer paragraph (with the "_GoBack" bookmark restored right
#    after "first finger!" and a trailing space run).
# ------------------------------------------------------------------
$ns = "xmlns:w=`"http://schemas.openxmlformats.org/wordprocessingml/2006/main`""

$pPrCommon = "<w:pPr><w:widowControl w:val=`"0`"/><w:tabs><w:tab w:val=`"left`" w:pos=`"220`"/><w:tab w:val=`"left`" w:pos=`"720`"/></w:tabs><w:autoSpaceDE w:val=`"0`"/><w:autoSpaceDN w:val=`"0`"/><w:adjustRightInd w:val=`"0`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr>"
$pPrIndented = "<w:pPr><w:widowControl w:val=`"0`"/><w:tabs><w:tab w:val=`"left`" w:pos=`"220`"/><w:tab w:val=`"left`" w:pos=`"720`"/></w:tabs><w:autoSpaceDE w:val=`"0`"/><w:autoSpaceDN w:val=`"0`"/><w:adjustRightInd w:val=`"0`"/><w:ind w:left=`"220`"/><w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr></w:pPr>"
$rPr = "<w:rPr><w:rFonts w:ascii=`"Arial`" w:hAnsi=`"Arial`" w:cs=`"Arial`"/></w:rPr>"

$blankPara = "<w:p/>"
$blankFormattedPara = "<w:p>$pPrCommon</w:p>"
$headingPara = "<w:p>$pPrCommon<w:r>$rPr<w:t>3.</w:t></w:r><w:r>$rPr<w:t xml:space=`"preserve`"> Identify potential solutions</w:t></w:r></w:p>"
$answerPara = "<w:p>$pPrIndented" + `
    "<w:r>$rPr<w:t xml:space=`"preserve`">The solution is </w:t></w:r>" + `
    "<w:r>$rPr<w:t>as follows</w:t></w:r>" + `
    "<w:r>$rPr<w:t xml:space=`"preserve`">. On every tenth count she is on her first finger. This means that since our problem is based on increments of 10 the answer is and always will be she is on her </w:t></w:r>" + `
    "<w:r>$rPr<w:t>first finger!</w:t></w:r>" + `
    "<w:bookmarkStart w:id=`"0`" w:name=`"_GoBack`"/><w:bookmarkEnd w:id=`"0`"/>" + `
    "<w:r>$rPr<w:t xml:space=`"preserve`"> </w:t></w:r>" + `
    "</w:p>"

$newXml = "<w:p $ns>".Replace("<w:p $ns>", "") # no-op placeholder to keep style consistent
$fullXml = $blankPara + $blankFormattedPara + $headingPara + $answerPara
$fullXml = $fullXml.Replace("<w:p>", "<w:p $ns>").Replace("<w:p/>", "<w:p $ns/>")

$insertAt.InsertXML($fullXml) | Out-Null
